$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 173
    "B3" = 194
    "B4" = 199
    "B5" = 220
    "B6" = 225
    "B7" = 230
    "B8" = 241
    "B9" = 246
    "B10" = 251
    "B11" = 257
    "B12" = 262
    "B13" = 267
    "B14" = 272
    "B15" = 288
    "B16" = 304
    "B17" = 325
    "B18" = 330
    "B19" = 335
    "B20" = 340
    "B21" = 346
    "B22" = 356
    "B23" = 361
    "B24" = 367
    "B25" = 372
    "B26" = 377
    "B27" = 382
    "B28" = 388
    "B29" = 393
    "B30" = 403
    "B31" = 409
    "B32" = 414
    "B33" = 419
    "B34" = 424
    "B35" = 445
    "B36" = 451
    "B37" = 456
    "B38" = 472
    "B39" = 477
    "B40" = 482
    "B41" = 487
    "B42" = 493
    "B43" = 519
    "B44" = 524
    "B45" = 540
    "B46" = 545
    "B47" = 566
    "B48" = 571
    "B49" = 577
    "B50" = 587
    "B51" = 592
    "B52" = 597
    "B53" = 613
    "B54" = 618
    "B55" = 634
    "B56" = 639
    "B57" = 645
    "B58" = 655
    "B59" = 660
    "B60" = 666
    "B61" = 750
    "B62" = 755
    "B63" = 760
    "B64" = 765
    "B65" = 771
    "B66" = 776
    "B67" = 781
    "B68" = 786
    "B69" = 792
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
